# Update LR-pair edge-weight metrics with newly recomputed TPM values.
# Source: NatmiData/natmiOut_TPM/YoungD0/LR-pairs_lrc2p/Inhba-Eng.xlsx
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 5).Value = 1
$ws.Cells.Item(2, 6).Value = 0.3333333333333333
$ws.Cells.Item(2, 7).Value = 0.05045533333333333
$ws.Cells.Item(2, 8).Value = 0.151366
$ws.Cells.Item(2, 9).Value = 0.004442474524580737
$ws.Cells.Item(2, 10).Value = 0.004442474524580737
$ws.Cells.Item(2, 13).Value = 135.0916853333333
$ws.Cells.Item(2, 14).Value = 405.2750559999999
$ws.Cells.Item(2, 15).Value = 0.7123704212620513
$ws.Cells.Item(2, 16).Value = 0.7123704212620514
$ws.Cells.Item(2, 17).Value = 6.81609601405511
$ws.Cells.Item(2, 18).Value = 61.34486412649599
$ws.Cells.Item(2, 19).Value = 0.003164687448521511
$ws.Cells.Item(2, 20).Value = 0.003164687448521511
$ws.Cells.Item(3, 5).Value = 1
$ws.Cells.Item(3, 6).Value = 0.3333333333333333
$ws.Cells.Item(3, 7).Value = 0.05045533333333333
$ws.Cells.Item(3, 8).Value = 0.151366
$ws.Cells.Item(3, 9).Value = 0.004442474524580737
$ws.Cells.Item(3, 10).Value = 0.004442474524580737
$ws.Cells.Item(3, 15).Value = 0.2125756143240238
$ws.Cells.Item(3, 16).Value = 0.2125756143240238
$ws.Cells.Item(3, 17).Value = 2.033964008376889
$ws.Cells.Item(3, 18).Value = 18.305676075392
$ws.Cells.Item(3, 19).Value = 0.0009443617511815757
$ws.Cells.Item(3, 20).Value = 0.0009443617511815758
$ws.Cells.Item(4, 5).Value = 1
$ws.Cells.Item(4, 6).Value = 0.3333333333333333
$ws.Cells.Item(4, 7).Value = 0.05045533333333333
$ws.Cells.Item(4, 8).Value = 0.151366
$ws.Cells.Item(4, 9).Value = 0.004442474524580737
$ws.Cells.Item(4, 10).Value = 0.004442474524580737
$ws.Cells.Item(4, 13).Value = 14.23299766666667
$ws.Cells.Item(4, 14).Value = 42.698993
$ws.Cells.Item(4, 15).Value = 0.07505396441392481
$ws.Cells.Item(4, 16).Value = 0.07505396441392483
$ws.Cells.Item(4, 17).Value = 0.7181306416042222
$ws.Cells.Item(4, 18).Value = 6.463175774438001
$ws.Cells.Item(4, 19).Value = 0.0003334253248776502
$ws.Cells.Item(4, 20).Value = 0.0003334253248776503
$ws.Cells.Item(5, 9).Value = 0.7425623198471305
$ws.Cells.Item(5, 10).Value = 0.7425623198471305
$ws.Cells.Item(5, 13).Value = 135.0916853333333
$ws.Cells.Item(5, 14).Value = 405.2750559999999
$ws.Cells.Item(5, 15).Value = 0.7123704212620513
$ws.Cells.Item(5, 16).Value = 0.7123704212620514
$ws.Cells.Item(5, 17).Value = 1139.314596064051
$ws.Cells.Item(5, 18).Value = 10253.83136457646
$ws.Cells.Item(5, 19).Value = 0.5289794326028264
$ws.Cells.Item(5, 20).Value = 0.5289794326028265
$ws.Cells.Item(6, 9).Value = 0.7425623198471305
$ws.Cells.Item(6, 10).Value = 0.7425623198471305
$ws.Cells.Item(6, 15).Value = 0.2125756143240238
$ws.Cells.Item(6, 16).Value = 0.2125756143240238
$ws.Cells.Item(6, 19).Value = 0.157850641315376
$ws.Cells.Item(6, 20).Value = 0.157850641315376
$ws.Cells.Item(7, 9).Value = 0.7425623198471305
$ws.Cells.Item(7, 10).Value = 0.7425623198471305
$ws.Cells.Item(7, 13).Value = 14.23299766666667
$ws.Cells.Item(7, 14).Value = 42.698993
$ws.Cells.Item(7, 15).Value = 0.07505396441392481
$ws.Cells.Item(7, 16).Value = 0.07505396441392483
$ws.Cells.Item(7, 17).Value = 120.0359736971741
$ws.Cells.Item(7, 18).Value = 1080.323763274567
$ws.Cells.Item(7, 19).Value = 0.05573224592892799
$ws.Cells.Item(7, 20).Value = 0.055732245928928
$ws.Cells.Item(8, 7).Value = 2.873389
$ws.Cells.Item(8, 8).Value = 8.620167
$ws.Cells.Item(8, 9).Value = 0.2529952056282888
$ws.Cells.Item(8, 10).Value = 0.2529952056282888
$ws.Cells.Item(8, 13).Value = 135.0916853333333
$ws.Cells.Item(8, 14).Value = 405.2750559999999
$ws.Cells.Item(8, 15).Value = 0.7123704212620513
$ws.Cells.Item(8, 16).Value = 0.7123704212620514
$ws.Cells.Item(8, 17).Value = 388.1709626282613
$ws.Cells.Item(8, 18).Value = 3493.538663654352
$ws.Cells.Item(8, 19).Value = 0.1802263012107034
$ws.Cells.Item(8, 20).Value = 0.1802263012107034
$ws.Cells.Item(9, 7).Value = 2.873389
$ws.Cells.Item(9, 8).Value = 8.620167
$ws.Cells.Item(9, 9).Value = 0.2529952056282888
$ws.Cells.Item(9, 10).Value = 0.2529952056282888
$ws.Cells.Item(9, 15).Value = 0.2125756143240238
$ws.Cells.Item(9, 16).Value = 0.2125756143240238
$ws.Cells.Item(9, 17).Value = 115.8325477597227
$ws.Cells.Item(9, 18).Value = 1042.492929837504
$ws.Cells.Item(9, 19).Value = 0.05378061125746621
$ws.Cells.Item(9, 20).Value = 0.05378061125746622
$ws.Cells.Item(10, 7).Value = 2.873389
$ws.Cells.Item(10, 8).Value = 8.620167
$ws.Cells.Item(10, 9).Value = 0.2529952056282888
$ws.Cells.Item(10, 10).Value = 0.2529952056282888
$ws.Cells.Item(10, 15).Value = 0.07505396441392481
$ws.Cells.Item(10, 16).Value = 0.07505396441392483
$ws.Cells.Item(10, 17).Value = 40.89693893242567
$ws.Cells.Item(10, 18).Value = 368.072450391831
$ws.Cells.Item(10, 19).Value = 0.01898829316011918
$ws.Cells.Item(10, 20).Value = 0.01898829316011918
